# DictionnaireDeDonnées.xlsx - "Add files via upload" edit
#
# Semantic change on sheet "Feuil1":
#  - The admin-related block (rows "IdAdmin" / "NomAdmin" / "PwdAdmin",
#    Excel rows 7-9) is removed entirely.
#  - The "PostUser" row (Excel row 6) is turned into a "Statut" row:
#    column A becomes "Statut", columns B and E are cleared (column C,
#    the Nature "A", stays as-is).
#  - The last row ("ListeAttente", now Excel row 20 after the delete)
#    switches its highlight color from yellow to orange.
#  - The active selection moves to the whole of row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the 3 "admin" rows (IdAdmin/NomAdmin/PwdAdmin), originally rows 7-9.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# Turn the "PostUser" row into the new "Statut" row.
$ws.Range("A6").Value = "Statut"
$ws.Range("B6").ClearContents()
$ws.Range("E6").ClearContents()

# The trailing "ListeAttente" row (now row 20) gets an orange fill
# (was yellow) to match its new position/meaning.
$ws.Range("A20:E20").Interior.Color = 49407

# Match the author's final selection: the whole 6th row.
$ws.Rows.Item(6).Select() | Out-Null
